{"js": "const body = context.document.body;\nconst results = body.search(\"Story Cards of your project\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text 'Story Cards of your project' not found\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"Project Story Cards of your project\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Story Cards of your project\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Project Story Cards of your project\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
